# Auto update Excel log
# Appends new sensor-log rows (2026-01-28, ~16:19-16:20) to the PIR,
# Humidity and Temperature sheets.
#
# Note: values like "2026-01-28" or "87.5%" would otherwise be
# auto-converted by Excel into a date serial / percentage number (with a
# matching number format). To keep them as literal text - matching the
# rest of the log - they are entered with a leading apostrophe (forces
# text entry) and the cell formatting is cleared right after so no stray
# number format sticks around on the new cells.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# PIR sheet: rows 107-119
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PIR")
$times = @("16:19:18","16:19:20","16:19:25","16:19:30","16:19:35","16:19:41","16:19:46","16:19:51","16:19:56","16:20:01","16:20:06","16:20:11","16:20:16")
$firstRow = 107
$r = $firstRow
foreach ($t in $times) {
    $ws.Cells.Item($r, 1).Value = "'2026-01-28"
    $ws.Cells.Item($r, 2).Value = $t
    $ws.Cells.Item($r, 3).Value = "16:00"
    $ws.Cells.Item($r, 4).Value = "Bathroom"
    $ws.Cells.Item($r, 5).Value = "No Motion"
    $ws.Cells.Item($r, 6).Value = "Inactive"
    $r = $r + 1
}
$lastRow = $r - 1
$ws.Range("A$($firstRow):F$($lastRow)").ClearFormats()

# ----------------------------------------------------------------------
# Humidity sheet: rows 108-118
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Humidity")
$times = @("16:19:18","16:19:20","16:19:24","16:19:28","16:19:36","16:19:40","16:19:48","16:20:00","16:20:04","16:20:08","16:20:12")
$values = @("87.5%","88.3%","88.4%","87.5%","87.4%","88.4%","87.4%","87.4%","88.3%","87.4%","88.3%")
$firstRow = 108
$r = $firstRow
for ($i = 0; $i -lt $times.Count; $i++) {
    $ws.Cells.Item($r, 1).Value = "'2026-01-28"
    $ws.Cells.Item($r, 2).Value = $times[$i]
    $ws.Cells.Item($r, 3).Value = "16:00"
    $ws.Cells.Item($r, 4).Value = "Bathroom"
    $ws.Cells.Item($r, 5).Value = "'" + $values[$i]
    $ws.Cells.Item($r, 6).Value = "Active"
    $r = $r + 1
}
$lastRow = $r - 1
$ws.Range("A$($firstRow):F$($lastRow)").ClearFormats()

# ----------------------------------------------------------------------
# Temperature sheet: rows 108-118
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Temperature")
$times = @("16:19:18","16:19:20","16:19:24","16:19:28","16:19:36","16:19:40","16:19:48","16:20:00","16:20:04","16:20:08","16:20:12")
$values = @("22.8C","22.7C","22.8C","22.8C","22.7C","22.7C","22.8C","22.8C","22.7C","22.7C","22.7C")
$firstRow = 108
$r = $firstRow
for ($i = 0; $i -lt $times.Count; $i++) {
    $ws.Cells.Item($r, 1).Value = "'2026-01-28"
    $ws.Cells.Item($r, 2).Value = $times[$i]
    $ws.Cells.Item($r, 3).Value = "16:00"
    $ws.Cells.Item($r, 4).Value = "Bathroom"
    $ws.Cells.Item($r, 5).Value = $values[$i]
    $ws.Cells.Item($r, 6).Value = "Active"
    $r = $r + 1
}
$lastRow = $r - 1
$ws.Range("A$($firstRow):F$($lastRow)").ClearFormats()
